$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.251.88"
$ws.Range("E2").Value = "  -0.49%  "
$ws.Range("D3").Value = "1.589.35"
$ws.Range("E3").Value = "  -0.31%  "
$ws.Range("D5").Value = "'212.16"
$ws.Range("E5").Value = "  +0.91%  "
$ws.Range("D6").Value = "'0.506"
$ws.Range("E6").Value = "  +0.08%  "
$ws.Range("D7").Value = "'0.999"
$ws.Range("E7").Value = "  -0.05%  "
$ws.Range("D8").Value = "'0.245"
$ws.Range("E8").Value = "  -0.37%  "
$ws.Range("E9").Value = "  -1.11%  "
$ws.Range("D10").Value = "'19.24"
$ws.Range("E10").Value = "  -1.63%  "
$ws.Range("E11").Value = "  +0.30%  "
$ws.Range("D12").Value = "1.812.91"
$ws.Range("E12").Value = "  -0.09%  "
$ws.Range("D13").Value = "1.575.26"
$ws.Range("E13").Value = "  +0.36%  "
$ws.Range("E14").Value = "  -1.49%  "
$ws.Range("E15").Value = "  -0.63%  "
$ws.Range("D16").Value = "'63.91"
$ws.Range("E16").Value = "  -0.90%  "
$ws.Range("D17").Value = "26.240.97"
$ws.Range("E17").Value = "  -0.53%  "
$ws.Range("D18").Value = "0.0₃0726"
$ws.Range("E18").Value = "  -0.31%  "
$ws.Range("D19").Value = "'7.44"
$ws.Range("E19").Value = "  -0.53%  "
$ws.Range("D20").Value = "'214.38"
$ws.Range("D21").Value = "'0.999"
$ws.Range("E21").Value = "  -0.11%  "
$ws.Range("E22").Value = "  -0.48%  "
$ws.Range("D23").Value = "'8.99"
$ws.Range("E23").Value = "  +0.62%  "
$ws.Range("D24").Value = "'2.13"
$ws.Range("E24").Value = "  -1.21%  "
$ws.Range("D25").Value = "'144.66"
$ws.Range("E25").Value = "  -0.08%  "
$ws.Range("E26").Value = "  -0.06%  "
$ws.Range("E27").Value = "  -1.00%  "
$ws.Range("E28").Value = "  -0.58%  "
$ws.Range("D29").Value = "'15.11"
$ws.Range("E29").Value = "  -1.08%  "
$ws.Range("E30").Value = "  -1.96%  "
$ws.Range("E31").Value = "  +0.57%  "
$ws.Range("E32").Value = "  -0.89%  "
$ws.Range("D33").Value = "1.417.52"
$ws.Range("E33").Value = "  +7.87%  "
$ws.Range("E35").Value = "  -0.53%  "
$ws.Range("D36").Value = "'0.592"
$ws.Range("E36").Value = "  -4.08%  "
$ws.Range("D37").Value = "'1.46"
$ws.Range("E37").Value = "  -1.31%  "
$ws.Range("E38").Value = "  -1.35%  "
$ws.Range("D39").Value = "'5.92"
$ws.Range("E39").Value = "  +5.17%  "
$ws.Range("E40").Value = "  +1.33%  "
$ws.Range("E41").Value = "  -0.09%  "
$ws.Range("D42").Value = "'0.943"
$ws.Range("E42").Value = "  -14.24%  "
$ws.Range("D43").Value = "'0.765"
$ws.Range("E43").Value = "  +0.20%  "
$ws.Range("E44").Value = "  -0.27%  "
$ws.Range("D45").Value = "1.723.77"
$ws.Range("E45").Value = "  -0.26%  "
$ws.Range("D46").Value = "'61.21"
$ws.Range("E46").Value = "  -2.11%  "
$ws.Range("D47").Value = "'85.90"
$ws.Range("E47").Value = "  -2.39%  "
$ws.Range("E48").Value = "  -0.13%  "
$ws.Range("E49").Value = "  -0.65%  "
$ws.Range("D50").Value = "'0.0969"
$ws.Range("E50").Value = "  -1.47%  "
$ws.Range("E51").Value = "  -0.02%  "
